$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X3").Value = -0.069999999999993179
$ws.Range("Y3").Value = "Down"

$ws.Range("A4").Value = 42641.890092592592
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = -11
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = -16
$ws.Range("E4").Value = 13848
$ws.Range("F4").Value = 723
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 26
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 18805
$ws.Range("L4").Value = 144
$ws.Range("M4").Value = 111
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 11
$ws.Range("P4").Value = "Named"
$ws.Range("Q4").Value = 64.728146835133757
$ws.Range("R4").Value = -32.1
$ws.Range("S4").Value = -0.0755
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -0.0025000000000000001
$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("U4").Value = 6.79
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 0
